$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the summary report title (A1) and the pay-period start/end date labels (G/H columns)
$ws.Range("A1").Value = "Summary report for 12/1/2019 through 12/14/2019"

# Force these as plain text (not auto-converted to date serials) while
# keeping the cells' original (unstyled/default) formatting.
$ws.Range("G5:G24").NumberFormat = "@"
$ws.Range("G5:G24").Value = "12/1/2019"
$ws.Range("H5:H24").NumberFormat = "@"
$ws.Range("H5:H24").Value = "12/14/2019"
$ws.Range("G5:H24").ClearFormats()

# Update the "hours" column (C) values for the new pay period
$ws.Range("C5").Value = 28.75
$ws.Range("C6").Value = 77.34999999999999
$ws.Range("C7").Value = 107
$ws.Range("C8").Value = 90.88
$ws.Range("C10").Value = 193.17
$ws.Range("C11").Value = 64.59999999999999
$ws.Range("C13").Value = 70.98
$ws.Range("C14").Value = 121.77
$ws.Range("C15").Value = 262.87
$ws.Range("C16").Value = 147.6
$ws.Range("C18").Value = 108.53
$ws.Range("C19").Value = 108.32
$ws.Range("C20").Value = 177.88
$ws.Range("C21").Value = 218.3
$ws.Range("C22").Value = 273.52
$ws.Range("C23").Value = 174.52

$wb.Save()
